# Generate Report for Handoff
# Adds a new row (for e735ad7f-0cf4-43d7-8d7c-0981e3dfd96e.md) to each of the
# three report sheets: Overview, zh-cn, de-de - mirroring the existing row
# for 6f3e377e-a28b-497a-bd8d-f1cab8d9e27d.md.

$wb = $excel.ActiveWorkbook

$commitHash = "5ebe81fa979625fe977b0b98d5bf73117309a067"
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/"
$newFile = "e735ad7f-0cf4-43d7-8d7c-0981e3dfd96e.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A..G
#   A File Name, B Path And Name, C Extension, D Publish URL,
#   E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range()

$rngOverview.Cells.Item(1,1).Value = "'" + $newFile
$rngOverview.Cells.Item(1,2).Value = "'e2e\" + $newFile
$rngOverview.Cells.Item(1,3).Value = "'.md"
$rngOverview.Cells.Item(1,4).Value = "'"
$rngOverview.Cells.Item(1,5).Value = "'Ready for handoff"
$rngOverview.Cells.Item(1,6).Value = "'Ready for handoff"
$rngOverview.Cells.Item(1,7).Value = "'2016-09-06 15:23:57"

$wsOverview.Hyperlinks.Add($rngOverview.Cells.Item(1,2), ($repoBase + $newFile), [Type]::Missing, [Type]::Missing, ("e2e\" + $newFile))

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A..P
#   A Source File Name, B File Extension, C Status, D Source Path,
#   E Priority, F Content Duplicate, G Latest Handoff File,
#   H Latest Handoff Datetime, I Latest Target File, J Latest Handback File,
#   K Latest Handback DateTime, L Reference Tokens, M To be localized,
#   N Dependency From, O Has metadata, P Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range()

$rngZhCn.Cells.Item(1,1).Value = "'" + $newFile
$rngZhCn.Cells.Item(1,2).Value = "'.md"
$rngZhCn.Cells.Item(1,3).Value = "'Ready for handoff"
$rngZhCn.Cells.Item(1,4).Value = "'e2e"
$rngZhCn.Cells.Item(1,5).Value = "'ht"
$rngZhCn.Cells.Item(1,6).Value = "'False"
$rngZhCn.Cells.Item(1,7).Value = "'e735ad7f-0cf4-43d7-8d7c-0981e3dfd96e.67ea4809027b6869f8ccea6411d9980c8a3feb52.zh-cn.xlf"
$rngZhCn.Cells.Item(1,8).Value = "'2016-09-06 15:23:52"
$rngZhCn.Cells.Item(1,9).Value = "'"
$rngZhCn.Cells.Item(1,10).Value = "'"
$rngZhCn.Cells.Item(1,11).Value = "'0001-01-01 00:00:00"
$rngZhCn.Cells.Item(1,12).Value = "'"
$rngZhCn.Cells.Item(1,13).Value = "'True"
$rngZhCn.Cells.Item(1,14).Value = "'"
$rngZhCn.Cells.Item(1,15).Value = "'False"
$rngZhCn.Cells.Item(1,16).Value = "'"

$wsZhCn.Hyperlinks.Add($rngZhCn.Cells.Item(1,1), ($repoBase + $newFile), [Type]::Missing, [Type]::Missing, $newFile)

# ---------------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn, different handoff file + datetime
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range()

$rngDeDe.Cells.Item(1,1).Value = "'" + $newFile
$rngDeDe.Cells.Item(1,2).Value = "'.md"
$rngDeDe.Cells.Item(1,3).Value = "'Ready for handoff"
$rngDeDe.Cells.Item(1,4).Value = "'e2e"
$rngDeDe.Cells.Item(1,5).Value = "'ht"
$rngDeDe.Cells.Item(1,6).Value = "'False"
$rngDeDe.Cells.Item(1,7).Value = "'e735ad7f-0cf4-43d7-8d7c-0981e3dfd96e.67ea4809027b6869f8ccea6411d9980c8a3feb52.de-de.xlf"
$rngDeDe.Cells.Item(1,8).Value = "'2016-09-06 15:23:57"
$rngDeDe.Cells.Item(1,9).Value = "'"
$rngDeDe.Cells.Item(1,10).Value = "'"
$rngDeDe.Cells.Item(1,11).Value = "'0001-01-01 00:00:00"
$rngDeDe.Cells.Item(1,12).Value = "'"
$rngDeDe.Cells.Item(1,13).Value = "'True"
$rngDeDe.Cells.Item(1,14).Value = "'"
$rngDeDe.Cells.Item(1,15).Value = "'False"
$rngDeDe.Cells.Item(1,16).Value = "'"

$wsDeDe.Hyperlinks.Add($rngDeDe.Cells.Item(1,1), ($repoBase + $newFile), [Type]::Missing, [Type]::Missing, $newFile)

Write-Output "Added handoff row for $newFile to Overview, zh-cn, de-de."
